$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C (rows 2-43) holds the "Förändrad" (last changed) date, stored as
# the date serial number 45728 (2025-03-12). Update it to 45729 (2025-03-13).
for ($row = 2; $row -le 43; $row++) {
    $ws.Cells.Item($row, 3).Value = 45729
}
